$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.042.20"
$ws.Range("E2").Value = "'  -1.99%  "
$ws.Range("D3").Value = "'2.684.12"
$ws.Range("E3").Value = "'  -2.80%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'549.42"
$ws.Range("E5").Value = "'  -4.70%  "
$ws.Range("D6").Value = "'157.44"
$ws.Range("E6").Value = "'  -1.94%  "
$ws.Range("E7").Value = "'  +0.08%  "
$ws.Range("E8").Value = "'  -2.40%  "
$ws.Range("E9").Value = "'  -4.40%  "
$ws.Range("E10").Value = "'  -2.51%  "
$ws.Range("E11").Value = "'  -4.88%  "
$ws.Range("D12").Value = "'5.12"
$ws.Range("E12").Value = "'  -12.65%  "
$ws.Range("D13").Value = "'3.158.37"
$ws.Range("E13").Value = "'  -2.92%  "
$ws.Range("D14").Value = "'26.03"
$ws.Range("E14").Value = "'  -4.93%  "
$ws.Range("D15").Value = "'62.906.07"
$ws.Range("E15").Value = "'  -1.64%  "
$ws.Range("E16").Value = "'  -3.96%  "
$ws.Range("D17").Value = "'2.684.45"
$ws.Range("E17").Value = "'  -3.07%  "
$ws.Range("D18").Value = "'11.94"
$ws.Range("E18").Value = "'  -2.04%  "
$ws.Range("D19").Value = "'4.56"
$ws.Range("E19").Value = "'  -5.93%  "
$ws.Range("D20").Value = "'342.84"
$ws.Range("E20").Value = "'  -4.31%  "
$ws.Range("D21").Value = "'6.31"
$ws.Range("E21").Value = "'  -5.27%  "
$ws.Range("D22").Value = "'0.998"
$ws.Range("E23").Value = "'  -4.95%  "
$ws.Range("D24").Value = "'63.38"
$ws.Range("E24").Value = "'  -2.82%  "
$ws.Range("E25").Value = "'  -2.15%  "
$ws.Range("E26").Value = "'  +0.27%  "
$ws.Range("D27").Value = "'8.16"
$ws.Range("E27").Value = "'  -5.50%  "
$ws.Range("D28").Value = "'0.0₃0850"
$ws.Range("E28").Value = "'  -8.02%  "
$ws.Range("E29").Value = "'  -2.49%  "
$ws.Range("E30").Value = "'  -3.59%  "
$ws.Range("D31").Value = "'7.00"
$ws.Range("E31").Value = "'  -5.18%  "
$ws.Range("D32").Value = "'165.16"
$ws.Range("E32").Value = "'  -1.90%  "
$ws.Range("E33").Value = "'  +0.03%  "
$ws.Range("D34").Value = "'4.80"
$ws.Range("E34").Value = "'  -3.92%  "
$ws.Range("D35").Value = "'19.53"
$ws.Range("E35").Value = "'  -3.33%  "
$ws.Range("E36").Value = "'  -6.08%  "
$ws.Range("D37").Value = "'1.77"
$ws.Range("D38").Value = "'339.64"
$ws.Range("E38").Value = "'  -3.02%  "
$ws.Range("D39").Value = "'6.15"
$ws.Range("E39").Value = "'  -3.52%  "
$ws.Range("D40").Value = "'0.932"
$ws.Range("E40").Value = "'  -7.62%  "
$ws.Range("D41").Value = "'38.14"
$ws.Range("E41").Value = "'  -2.97%  "
$ws.Range("E42").Value = "'  -6.22%  "
$ws.Range("D43").Value = "'20.37"
$ws.Range("E43").Value = "'  -5.57%  "
$ws.Range("E44").Value = "'  -8.04%  "
$ws.Range("E45").Value = "'  -2.40%  "
$ws.Range("D46").Value = "'0.0560"
$ws.Range("E46").Value = "'  -5.90%  "
$ws.Range("E47").Value = "'  -0.01%  "
$ws.Range("D48").Value = "'11.03"
$ws.Range("E48").Value = "'  +0.01%  "
$ws.Range("E49").Value = "'  -4.05%  "
$ws.Range("D50").Value = "'128.98"
$ws.Range("E50").Value = "'  -5.73%  "
$ws.Range("D51").Value = "'2.085.45"
$ws.Range("E51").Value = "'  -2.93%  "
